# Auto-generated Excel COM-interop script to update cryptos list values
# Updates Price (D) and Volume(1h) (E) columns, and swaps ARBITRUM/ImmutableX rows (37/38)
# D-column cells hold plain numeric-looking strings (e.g. "27.206.21"), so we force
# the Text number format before assigning them to avoid Excel auto-converting them to
# numbers (which would strip significant trailing zeros / renormalize the text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.206.21'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.630.77'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.11'
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.520'
$ws.Range('E6').Value = '  +1.66%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.256'
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0625'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.26'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0847'
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.649.97'
$ws.Range('E12').Value = '  -1.25%  '
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.177.94'
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.81'
$ws.Range('E16').Value = '  -4.87%  '
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '216.37'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.93'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('E21').Value = '  -1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.49'
$ws.Range('E22').Value = '  -0.88%  '
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '147.78'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.29'
$ws.Range('E26').Value = '  -3.75%  '
$ws.Range('E27').Value = '  -0.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.60'
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('E32').Value = '  -1.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.319.32'
$ws.Range('E33').Value = '  +4.84%  '
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.45'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.850'
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.541'
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('E40').Value = '  +1.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.800'
$ws.Range('E41').Value = '  -1.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '63.75'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.768.46'
$ws.Range('E43').Value = '  -1.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.20'
$ws.Range('E44').Value = '  -4.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.72'
$ws.Range('E45').Value = '  -1.16%  '
$ws.Range('E47').Value = '  +5.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.809'
$ws.Range('E48').Value = '  +21.09%  '
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.56'
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0956'
$ws.Range('E51').Value = '  -2.02%  '
